# Fruta / hortaliza, semanal
# Insert one new weekly data row at row 288 (pushing the existing rows
# 288-313 down to 289-314 unchanged) on the "Papa" price sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 288..313 down to 289..314 by inserting a blank row at 288.
$ws.Rows.Item(288).Insert()

# Populate the newly inserted row 288 with this week's record.
$ws.Range("A288").Value = 4
$ws.Range("B288").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C288").Value = "Los Lagos"
$ws.Range("D288").Value = 44578
$ws.Range("E288").Value = 10
$ws.Range("F288").Value = 100114001
$ws.Range("G288").Value = "Papa"
$ws.Range("H288").Value = "Patagonia"
$ws.Range("I288").Value = "1a nueva(o)"
$ws.Range("J288").Value = 250
$ws.Range("K288").Value = 9000
$ws.Range("L288").Value = 10000
$ws.Range("M288").Value = 9400
$ws.Range("N288").Value = "`$/saco 25 kilos"
$ws.Range("O288").Value = "Provincia de Llanquihue"
$ws.Range("P288").Value = 376
$ws.Range("Q288").Value = 25
$ws.Range("R288").Value = "Hortaliza"
